# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worksheet "Hoja1" has a table of worker account-statement rows
# (B:Tipo Doc, C:N Doc, D:Nombre, E:Periodo Mora, F:Valor Mora, G:Salario
# Basico) starting at row 16. This edit:
#   - moves the existing two workers (EFREN JULIO VILLAMIL / ANIBAL SENEN
#     PAJARO DE AVILA) down to the bottom of the table (rows 26-27) and
#     refreshes their Salario Basico values,
#   - inserts ten new rows (16-25) for MIGUEL ANGEL MARTINEZ VASCO, one per
#     period 1803..1812, each with Valor Mora 44000 / Salario Basico
#     1100000.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New block: MIGUEL ANGEL MARTINEZ VASCO, periods 1803-1812 (rows 16-25)
$periodos = @("1803","1804","1805","1806","1807","1808","1809","1810","1811","1812")
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 2).Value = "CC"
    $ws.Cells.Item($row, 3).Value = "1026156533"
    $ws.Cells.Item($row, 4).Value = "MIGUEL ANGEL MARTINEZ VASCO"
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = 44000
    $ws.Cells.Item($row, 7).Value = 1100000
}

# Existing worker 1, now at row 26, with updated Salario Basico
$ws.Cells.Item(26, 2).Value = "CC"
$ws.Cells.Item(26, 3).Value = "73137452"
$ws.Cells.Item(26, 4).Value = "EFREN JULIO VILLAMIL"
$ws.Cells.Item(26, 5).Value = "1911"
$ws.Cells.Item(26, 6).Value = 40000
$ws.Cells.Item(26, 7).Value = 1423500

# Existing worker 2, now at row 27, with updated Salario Basico
$ws.Cells.Item(27, 2).Value = "CC"
$ws.Cells.Item(27, 3).Value = "73131006"
$ws.Cells.Item(27, 4).Value = "ANIBAL SENEN PAJARO DE AVILA"
$ws.Cells.Item(27, 5).Value = "2203"
$ws.Cells.Item(27, 6).Value = 50000
$ws.Cells.Item(27, 7).Value = 2000000
